$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.855412382637946
$ws.Range("C2").Value = 5.72609263032619
$ws.Range("D2").Value = 10.99082523556015
$ws.Range("F2").Value = 31.73380328294286
$ws.Range("G2").Value = 31.07709316244517
$ws.Range("H2").Value = 14.9983615352189
$ws.Range("J2").Value = 11.51203972391945
$ws.Range("K2").Value = 8.081647105149639
$ws.Range("M2").Value = 15.46151087255969
$ws.Range("O2").Value = 23.11399910700465
$ws.Range("B3").Value = 8.554288718961949
$ws.Range("C3").Value = 5.633480733610848
$ws.Range("D3").Value = 10.94305937733673
$ws.Range("F3").Value = 31.81675099005536
$ws.Range("G3").Value = 31.19942635312756
$ws.Range("H3").Value = 15.04991577096695
$ws.Range("J3").Value = 11.53285902276234
$ws.Range("K3").Value = 7.870683070697826
$ws.Range("M3").Value = 15.36976518107338
$ws.Range("O3").Value = 23.20410715127628
$ws.Range("B4").Value = 8.364838524017589
$ws.Range("C4").Value = 5.575664925535744
$ws.Range("D4").Value = 10.9155521400319
$ws.Range("F4").Value = 31.87485089427859
$ws.Range("G4").Value = 31.28350742403019
$ws.Range("H4").Value = 15.08376519329008
$ws.Range("J4").Value = 11.54733060962598
$ws.Range("K4").Value = 7.738893443013051
$ws.Range("M4").Value = 15.31543884710581
$ws.Range("O4").Value = 23.26394056452186
$ws.Range("B5").Value = 8.286606199505792
$ws.Range("C5").Value = 5.551885881662214
$ws.Range("D5").Value = 10.90480926134118
$ws.Range("F5").Value = 31.90032655612521
$ws.Range("G5").Value = 31.32001771082777
$ws.Range("H5").Value = 15.09811140334428
$ws.Range("J5").Value = 11.55365253005139
$ws.Range("K5").Value = 7.68469487047622
$ws.Range("M5").Value = 15.29382212195273
$ws.Range("O5").Value = 23.28945534708104
$ws.Range("B6").Value = 8.273557194700015
$ws.Range("C6").Value = 5.547924800555048
$ws.Range("D6").Value = 10.90305383265924
$ws.Range("F6").Value = 31.90466534921616
$ws.Range("G6").Value = 31.32621564846514
$ws.Range("H6").Value = 15.10052694887663
$ws.Range("J6").Value = 11.55472792765199
$ws.Range("K6").Value = 7.675667702777876
$ws.Range("M6").Value = 15.2902646814999
$ws.Range("O6").Value = 23.29376040702282
$ws.Range("B7").Value = 8.363787463535518
$ws.Range("C7").Value = 5.575345090511504
$ws.Range("D7").Value = 10.91540535785786
$ws.Range("F7").Value = 31.87518718666181
$ws.Range("G7").Value = 31.28399072974638
$ws.Range("H7").Value = 15.08395643458099
$ws.Range("J7").Value = 11.54741414989458
$ws.Range("K7").Value = 7.738164397908238
$ws.Range("M7").Value = 15.31514518175161
$ws.Range("O7").Value = 23.26428008306238
$ws.Range("B8").Value = 8.752597874119605
$ws.Range("C8").Value = 5.694369802169692
$ws.Range("D8").Value = 10.9739828976016
$ws.Range("F8").Value = 31.76091363540717
$ws.Range("G8").Value = 31.11740731895641
$ws.Range("H8").Value = 15.01568210189969
$ws.Range("J8").Value = 11.51886782444856
$ws.Range("K8").Value = 8.009418338245844
$ws.Range("M8").Value = 15.42947122870657
$ws.Range("O8").Value = 23.14413235501863
$ws.Range("B9").Value = 9.474288189260786
$ws.Range("C9").Value = 5.919349680474982
$ws.Range("D9").Value = 11.10290201355982
$ws.Range("F9").Value = 31.59385192834093
$ws.Range("G9").Value = 30.86226137770822
$ws.Range("H9").Value = 14.89919538640014
$ws.Range("J9").Value = 11.47628236258265
$ws.Range("K9").Value = 8.520401232973231
$ws.Range("M9").Value = 15.6687890157839
$ws.Range("O9").Value = 22.94432906163782
$ws.Range("B10").Value = 9.974234479495033
$ws.Range("C10").Value = 6.078384197881974
$ws.Range("D10").Value = 11.20561238376841
$ws.Range("F10").Value = 31.50605185463576
$ws.Range("G10").Value = 30.71888825792101
$ws.Range("H10").Value = 14.82419647765323
$ws.Range("J10").Value = 11.45315420646487
$ws.Range("K10").Value = 8.879313778416662
$ws.Range("M10").Value = 15.85273653843988
$ws.Range("O10").Value = 22.8194259836811
$ws.Range("B11").Value = 10.19414615454704
$ws.Range("C11").Value = 6.149152809022108
$ws.Range("D11").Value = 11.2539424650381
$ws.Range("F11").Value = 31.47372630404013
$ws.Range("G11").Value = 30.66333097881952
$ws.Range("H11").Value = 14.7923710690242
$ws.Range("J11").Value = 11.44440223309384
$ws.Range("K11").Value = 9.038313991296205
$ws.Range("M11").Value = 15.93793765286404
$ws.Range("O11").Value = 22.76737184212975
$ws.Range("B12").Value = 10.27627285035704
$ws.Range("C12").Value = 6.17570690584512
$ws.Range("D12").Value = 11.27246313456239
$ws.Range("F12").Value = 31.46258230636978
$ws.Range("G12").Value = 30.64368964535692
$ws.Range("H12").Value = 14.78064890535934
$ws.Range("J12").Value = 11.44134221907172
$ws.Range("K12").Value = 9.097857692351647
$ws.Range("M12").Value = 15.97039879691992
$ws.Range("O12").Value = 22.74834668920817
$ws.Range("B13").Value = 10.25863750214436
$ws.Range("C13").Value = 6.169999148981372
$ws.Range("D13").Value = 11.26846482506139
$ws.Range("F13").Value = 31.46493355048414
$ws.Range("G13").Value = 30.64785751409067
$ws.Range("H13").Value = 14.78315883397907
$ws.Range("J13").Value = 11.4419899479043
$ws.Range("K13").Value = 9.085064301023104
$ws.Range("M13").Value = 15.9633992790736
$ws.Range("O13").Value = 22.75241353592033
$ws.Range("B14").Value = 10.20092615822208
$ws.Range("C14").Value = 6.151342427402959
$ws.Range("D14").Value = 11.25546184970119
$ws.Range("F14").Value = 31.47278748578463
$ws.Range("G14").Value = 30.66168703613183
$ws.Range("H14").Value = 14.79140007818138
$ws.Range("J14").Value = 11.44414539150614
$ws.Range("K14").Value = 9.043226321248302
$ws.Range("M14").Value = 15.94060442056515
$ws.Range("O14").Value = 22.76579286226647
$ws.Range("B15").Value = 10.16542476468596
$ws.Range("C15").Value = 6.139882304821798
$ws.Range("D15").Value = 11.247525327825
$ws.Range("F15").Value = 31.47774115837155
$ws.Range("G15").Value = 30.67034015842573
$ws.Range("H15").Value = 14.79649097507047
$ws.Range("J15").Value = 11.44549875536979
$ws.Range("K15").Value = 9.017511122555213
$ws.Range("M15").Value = 15.92666696428878
$ws.Range("O15").Value = 22.7740775466636
$ws.Range("B16").Value = 9.9597054624327
$ws.Range("C16").Value = 6.073726037897311
$ws.Range("D16").Value = 11.20248521826545
$ws.Range("F16").Value = 31.50831781677347
$ws.Range("G16").Value = 30.72271415789151
$ws.Range("H16").Value = 14.82632244975566
$ws.Range("J16").Value = 11.45376174954146
$ws.Range("K16").Value = 8.868832239542305
$ws.Range("M16").Value = 15.84719724137438
$ws.Range("O16").Value = 22.82292384899131
$ws.Range("B17").Value = 9.831526992175336
$ws.Range("C17").Value = 6.032724539652541
$ws.Range("D17").Value = 11.17525770122211
$ws.Range("F17").Value = 31.52902753840489
$ws.Range("G17").Value = 30.75732444160974
$ws.Range("H17").Value = 14.84521000792367
$ws.Range("J17").Value = 11.4592838001826
$ws.Range("K17").Value = 8.776490122943079
$ws.Range("M17").Value = 15.79881980675446
$ws.Range("O17").Value = 22.85411090118335
$ws.Range("B18").Value = 9.757098173841641
$ws.Range("C18").Value = 6.008994127468509
$ws.Range("D18").Value = 11.15974911300741
$ws.Range("F18").Value = 31.54165590593676
$ws.Range("G18").Value = 30.77814056100594
$ws.Range("H18").Value = 14.8562893714333
$ws.Range("J18").Value = 11.4626264748038
$ws.Range("K18").Value = 8.722978190811064
$ws.Range("M18").Value = 15.77113900065102
$ws.Range("O18").Value = 22.87249727796039
$ws.Range("B19").Value = 9.731779166256064
$ws.Range("C19").Value = 6.000934649186027
$ws.Range("D19").Value = 11.15452462784379
$ws.Range("F19").Value = 31.54605467939785
$ws.Range("G19").Value = 30.78534446546076
$ws.Range("H19").Value = 14.86007771196379
$ws.Range("J19").Value = 11.46378685610284
$ws.Range("K19").Value = 8.704793073226504
$ws.Range("M19").Value = 15.76179224637685
$ws.Range("O19").Value = 22.8787995600718
$ws.Range("B20").Value = 9.845245186133884
$ws.Range("C20").Value = 6.037104609040786
$ws.Range("D20").Value = 11.17814047484395
$ws.Range("F20").Value = 31.52674876376533
$ws.Range("G20").Value = 30.75354596428521
$ws.Range("H20").Value = 14.8431770680286
$ws.Range("J20").Value = 11.45867873400601
$ws.Range("K20").Value = 8.786361805646051
$ws.Range("M20").Value = 15.80395485572888
$ws.Range("O20").Value = 22.85074456713152
$ws.Range("B21").Value = 10.21790906773492
$ws.Range("C21").Value = 6.156829129025819
$ws.Range("D21").Value = 11.25927528864508
$ws.Range("F21").Value = 31.47045080933111
$ws.Range("G21").Value = 30.65758699989044
$ws.Range("H21").Value = 14.78897048558188
$ws.Range("J21").Value = 11.44350538953138
$ws.Range("K21").Value = 9.05553362309707
$ws.Range("M21").Value = 15.94729464099369
$ws.Range("O21").Value = 22.76184438372778
$ws.Range("B22").Value = 10.45474272855607
$ws.Range("C22").Value = 6.233643702890868
$ws.Range("D22").Value = 11.31357300078674
$ws.Range("F22").Value = 31.44005152037668
$ws.Range("G22").Value = 30.60301812638855
$ws.Range("H22").Value = 14.7554634451844
$ws.Range("J22").Value = 11.43507007355382
$ws.Range("K22").Value = 9.227548793614766
$ws.Range("M22").Value = 16.04211548634934
$ws.Range("O22").Value = 22.70774578368814
$ws.Range("B23").Value = 10.32897581887849
$ws.Range("C23").Value = 6.192782972980172
$ws.Range("D23").Value = 11.28448097702421
$ws.Range("F23").Value = 31.45569054812262
$ws.Range("G23").Value = 30.63139491454322
$ws.Range("H23").Value = 14.77317114235634
$ws.Range("J23").Value = 11.43943670817084
$ws.Range("K23").Value = 9.136114291898281
$ws.Range("M23").Value = 15.99141070530119
$ws.Range("O23").Value = 22.73625249623642
$ws.Range("B24").Value = 9.839045481932503
$ws.Range("C24").Value = 6.035124870790234
$ws.Range("D24").Value = 11.17683672072853
$ws.Range("F24").Value = 31.52777674872744
$ws.Range("G24").Value = 30.75525135429415
$ws.Range("H24").Value = 14.84409547281242
$ws.Range("J24").Value = 11.45895176118463
$ws.Range("K24").Value = 8.781900132126372
$ws.Range("M24").Value = 15.80163288720486
$ws.Range("O24").Value = 22.8522650648121
$ws.Range("B25").Value = 9.284015195879794
$ws.Range("C25").Value = 5.85950760583436
$ws.Range("D25").Value = 11.06657898181914
$ws.Range("F25").Value = 31.63292185868129
$ws.Range("G25").Value = 30.92357693930083
$ws.Range("H25").Value = 14.92884787347459
$ws.Range("J25").Value = 11.48636937475245
$ws.Range("K25").Value = 8.384813784377558
$ws.Range("M25").Value = 15.60253933081584
$ws.Range("O25").Value = 22.99454123569955
